$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text describing the population assumption (finite -> infinite)
$ws.Range("B21").Value = "A população avaliada é considerada infinita"

# Update the recalculated numeric results (population now treated as infinite)
$ws.Range("B2").Value  = 21.79714285714286
$ws.Range("B3").Value  = 21.720718880312205
$ws.Range("B4").Value  = 21.720718880312205
$ws.Range("B5").Value  = 217971.42857142858
$ws.Range("B6").Value  = 217971.42857142858
$ws.Range("B7").Value  = 217971.42857142858
$ws.Range("B8").Value  = 4.71386089662237
$ws.Range("B9").Value  = 1.0
$ws.Range("B10").Value = 0.0
$ws.Range("B11").Value = 1.0274869937237732
$ws.Range("B12").Value = 8.626819848816465
$ws.Range("B13").Value = 74.42202070393374
$ws.Range("B14").Value = 1.055729522271517
$ws.Range("B15").Value = 22.220484552705457
$ws.Range("B16").Value = 39.57775523772136
$ws.Range("B17").Value = 2.179714285714286
$ws.Range("B18").Value = 1.0
$ws.Range("B19").Value = 0.993
$ws.Range("B20").Value = 0.0
$ws.Range("B22").Value = 10000.0
$ws.Range("B23").Value = 1.0
